$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Fill in new journal entry on row 22 ---
$ws.Range("B22").Value = 7
$ws.Range("C22").Value = "Recherche de solutions pour comment mettre en ordre par le meilleur score au pire dans le tableau des scores"
$ws.Range("D22").Value = "1h30"
$ws.Range("F22").Value = "toujours entrain d'essayer d'en trouver une"

# Task/resolution cells use the wrapped, centered style (same as the comment column)
$ws.Range("C22").WrapText = $true
$ws.Range("C22").HorizontalAlignment = -4108
$ws.Range("C22").VerticalAlignment = -4108
$ws.Range("F22").WrapText = $true
$ws.Range("F22").HorizontalAlignment = -4108
$ws.Range("F22").VerticalAlignment = -4108

# Row 22 grows taller to fit the wrapped text, matching the other data rows
$ws.Rows.Item(22).RowHeight = 28.8

# --- Update the view state (scrolled/selected cell moved down-right) ---
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F22").Select()
